$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.7173212604331821, 0.4491229322558946, 1.4172469810181, 2.467447866022739)"
$ws.Range("C2").Value = "NCT(2.313023586362248, 1.1834717658602893, 2.614935816503569, 4.001514159526913)"
$ws.Range("D2").Value = "F(836.2363470994462, 17.635217929068403, -9.922188133077935, 13.097282845854316)"
$ws.Range("E2").Value = "NCT(2.9213136182917894, 1.8429487208744229, -0.01492931809437336, 4.5641672076154)"
